$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Populate new hour entries (F column for Blad1 period 2, B/C entries)
$ws.Range("E18").Value = 41175
$ws.Range("F18").Value = 7

$ws.Range("E19").Value = 41193
$ws.Range("F19").Value = 2

$ws.Range("B24").Value = 41175
$ws.Range("C24").Value = 7

# Update selection to reflect last edited cell
$ws.Range("S21").Select()
